$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14

$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"
$ws.Cells.Item($row, 4).Value = 44476
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(6, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = 300000000
$ws.Cells.Item($row, 7).Value = "Espárragos"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 60
$ws.Cells.Item($row, 11).Value = 2000
$ws.Cells.Item($row, 12).Value = 2000
$ws.Cells.Item($row, 13).Value = 2000
$ws.Cells.Item($row, 14).Value = "$/kilo"
$ws.Cells.Item($row, 15).Value = "Provincia de Linares"
$ws.Cells.Item($row, 16).Value = 2000
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
